$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column T as a copy of column S (preserves cell styles exactly)
$ws.Columns("S").Copy()
$ws.Columns("T").Insert(-4161)
$excel.CutCopyMode = $false

# Row 3 header
$ws.Range("T3").Value = 2023

$ws.Range("T4").Value = 52.734251206028382
$ws.Range("T5").Value = 44.646801162600475
$ws.Range("T6").Value = 60.998061560200554
$ws.Range("T7").Value = 41.931627189714625
$ws.Range("T8").Value = 38.177163051511151
$ws.Range("T9").Value = 45.607453560981966
$ws.Range("T10").Value = 50.172884880431361
$ws.Range("T11").Value = 44.112367891063748
$ws.Range("T12").Value = 56.155144351753421
$ws.Range("T13").Value = 37.12775271808399
$ws.Range("T14").Value = 26.579446704517768
$ws.Range("T15").Value = 47.691579663423148
$ws.Range("T16").Value = 49.966474107695483
$ws.Range("T17").Value = 44.339536521432947
$ws.Range("T18").Value = 55.468421253968863
$ws.Range("T19").Value = 45.077411133103766
$ws.Range("T20").Value = 42.341975649266388
$ws.Range("T21").Value = 47.772457765110225
$ws.Range("T22").Value = 46.127136558116561
$ws.Range("T23").Value = 38.861148383596195
$ws.Range("T24").Value = 53.259250196123595
$ws.Range("T25").Value = 77.010722119335071
$ws.Range("T26").Value = 59.189709472566221
$ws.Range("T27").Value = 95.060584781341987
$ws.Range("T28").Value = 61.976853895626128
$ws.Range("T29").Value = 48.390234028455353
$ws.Range("T30").Value = 78.407224173903401
$ws.Range("T31").Value = 39.559787476030614
$ws.Range("T32").Value = 41.059215262778977
$ws.Range("T33").Value = 38.128175110356899

$ws.Range("A1").Select()